$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "Instituição" header from E2 to D2 and re-merge D2:F2 ---
$ws.Range("E2:F2").UnMerge()
$ws.Range("E2").Copy()
$ws.Range("D2").PasteSpecial(-4104) | Out-Null
$ws.Range("E2").ClearContents()
$ws.Range("D2:F2").Merge()

# --- Add a new "idTipoUsuario" foreign-key column (D) to the Instituição table ---
$ws.Range("E3").Copy()
$ws.Range("D3").PasteSpecial(-4104) | Out-Null
$ws.Range("D3").Value = "idTipoUsuario"

$ws.Range("E4").Copy()
$ws.Range("D4").PasteSpecial(-4104) | Out-Null
$ws.Range("D4").Value = 1

$ws.Range("E5").Copy()
$ws.Range("D5").PasteSpecial(-4104) | Out-Null
$ws.Range("D5").Value = 1

# --- Fix mislabeled column header in the Autor lookup table ---
$ws.Range("C23").Value = "idUsuario"

# --- Restore the view state to match what was last selected/scrolled ---
$ws.Range("E8").Select()
$ws.Application.ActiveWindow.ScrollColumn = 5
